# Auto-generated edit script applying numeric updates to the Brynhildr Profits workbook.
# Each assignment mirrors one cell-level change from the source diff; ClearContents()
# is used where a cell is removed entirely (no replacement value).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4033.5
$ws.Range("I12").Value = 176.42857
$ws.Range("J12").Value = 13033.333
$ws.Range("K12").Value = 176.42857
$ws.Range("L12").Value = 13033.333
$ws.Range("M12").Value = -6.428570000000008
$ws.Range("N12").Value = -13373.333
$ws.Range("H129").Value = 1609.6
$ws.Range("J129").Value = 2200
$ws.Range("L129").Value = 6600
$ws.Range("N129").Value = -16600
$ws.Range("H132").Value = 965.3333
$ws.Range("I132").Value = 886.8461
$ws.Range("K132").Value = 2660.5383
$ws.Range("M132").Value = -130.5383000000002
$ws.Range("H137").Value = 4612.1943
$ws.Range("I137").Value = 2102.4092
$ws.Range("J137").Value = 8556.143
$ws.Range("K137").Value = 6307.2276
$ws.Range("L137").Value = 25668.429
$ws.Range("M137").Value = -3757.2276
$ws.Range("N137").Value = -30768.429
$ws.Range("H138").Value = 2943.6858
$ws.Range("I138").Value = 2767.0667
$ws.Range("J138").Value = 4003.4
$ws.Range("K138").Value = 8301.2001
$ws.Range("L138").Value = 12010.2
$ws.Range("M138").Value = -3161.2001
$ws.Range("N138").Value = -22290.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8825.833000000001
$ws.Range("I61").Value = 10690.154
$ws.Range("J61").Value = 3978.6
$ws.Range("K61").Value = 10690.154
$ws.Range("L61").Value = 3978.6
$ws.Range("M61").Value = -10478.154
$ws.Range("N61").Value = -4402.6
$ws.Range("H74").Value = 5566.0557
$ws.Range("I74").Value = 4081.775
$ws.Range("K74").Value = 4081.775
$ws.Range("M74").Value = -3207.775
$ws.Range("H77").Value = 5566.0557
$ws.Range("I77").Value = 4081.775
$ws.Range("K77").Value = 20408.875
$ws.Range("M77").Value = -16040.875
$ws.Range("H97").Value = 38462612
$ws.Range("I97").Value = 861.5238000000001
$ws.Range("K97").Value = 861.5238000000001
$ws.Range("M97").Value = -365.5238000000001
$ws.Range("H132").Value = 5942.114
$ws.Range("I132").Value = 4290.36
$ws.Range("K132").Value = 12871.08
$ws.Range("M132").Value = -10341.08
$ws.Range("H133").Value = 82500
$ws.Range("J133").Value = 82500
$ws.Range("L133").Value = 82500
$ws.Range("N133").Value = -87560
$ws.Range("H136").Value = 8825.833000000001
$ws.Range("I136").Value = 10690.154
$ws.Range("J136").Value = 3978.6
$ws.Range("K136").Value = 32070.462
$ws.Range("L136").Value = 11935.8
$ws.Range("M136").Value = -29520.462
$ws.Range("N136").Value = -17035.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 793
$ws.Range("I11").Value = 91.333336
$ws.Range("K11").Value = 91.333336
$ws.Range("M11").Value = 48.666664
$ws.Range("H81").Value = 33572.715
$ws.Range("J81").Value = 33572.715
$ws.Range("L81").Value = 33572.715
$ws.Range("N81").Value = -35694.715
$ws.Range("H84").Value = 33572.715
$ws.Range("J84").Value = 33572.715
$ws.Range("L84").Value = 100718.145
$ws.Range("N84").Value = -111326.145
$ws.Range("H135").Value = 62051.168
$ws.Range("J135").Value = 62051.168
$ws.Range("L135").Value = 62051.168
$ws.Range("N135").Value = -72191.16800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 38150
$ws.Range("J50").Value = 38150
$ws.Range("L50").Value = 38150
$ws.Range("N50").Value = -39400
$ws.Range("H51").Value = 26500
$ws.Range("J51").Value = 26500
$ws.Range("L51").Value = 26500
$ws.Range("N51").Value = -27972
$ws.Range("H59").Value = 29427.857
$ws.Range("I59").Value = 25498.75
$ws.Range("J59").Value = 34666.668
$ws.Range("K59").Value = 25498.75
$ws.Range("L59").Value = 34666.668
$ws.Range("M59").Value = -24353.75
$ws.Range("N59").Value = -36956.668
$ws.Range("H60").Value = 20975
$ws.Range("J60").Value = 22966.666
$ws.Range("L60").Value = 22966.666
$ws.Range("N60").Value = -23988.666
$ws.Range("H61").Value = 26500
$ws.Range("J61").Value = 26500
$ws.Range("L61").Value = 26500
$ws.Range("N61").Value = -27196
$ws.Range("H124").Value = 31841.143
$ws.Range("J124").Value = 31841.143
$ws.Range("L124").Value = 31841.143
$ws.Range("N124").Value = -36751.143

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 7295.75
$ws.Range("J127").Value = 7295.75
$ws.Range("L127").Value = 21887.25
$ws.Range("N127").Value = -31807.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4221.625
$ws.Range("I21").Value = 4221.625
$ws.Range("K21").Value = 4221.625
$ws.Range("M21").Value = -4048.625
$ws.Range("H25").Value = 5009.3335
$ws.Range("I25").Value = 5009.3335
$ws.Range("K25").Value = 5009.3335
$ws.Range("M25").Value = -4480.3335
$ws.Range("H30").Value = 4221.625
$ws.Range("I30").Value = 4221.625
$ws.Range("K30").Value = 4221.625
$ws.Range("M30").Value = -4116.625
$ws.Range("H94").Value = 45000
$ws.Range("J94").Value = 45000
$ws.Range("L94").Value = 45000
$ws.Range("N94").Value = -46352
$ws.Range("H107").Value = 1107.5
$ws.Range("I107").Value = 1107.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1107.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 812.5
$ws.Range("N107").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 12110.739
$ws.Range("I132").Value = 8936.223
$ws.Range("K132").Value = 26808.669
$ws.Range("M132").Value = -24278.669
$ws.Range("H138").Value = 54950
$ws.Range("J138").Value = 54950
$ws.Range("L138").Value = 54950
$ws.Range("N138").Value = -65230

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3452.9473
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 3753.2942
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 3753.2942
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -4343.2942
$ws.Range("H27").Value = 3452.9473
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 3753.2942
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 3753.2942
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -3967.2942
$ws.Range("H62").Value = 25000
$ws.Range("I62").Value = 25000
$ws.Range("K62").Value = 25000
$ws.Range("M62").Value = -24376
$ws.Range("H65").Value = 25000
$ws.Range("I65").Value = 25000
$ws.Range("K65").Value = 75000
$ws.Range("M65").Value = -71880
$ws.Range("H122").Value = 5158.4165
$ws.Range("I122").Value = 3326
$ws.Range("J122").Value = 5769.222
$ws.Range("K122").Value = 9978
$ws.Range("L122").Value = 17307.666
$ws.Range("M122").Value = -7528
$ws.Range("N122").Value = -22207.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 215409
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H132").Value = 2775.348
$ws.Range("I132").Value = 2712.524
$ws.Range("K132").Value = 8137.572
$ws.Range("M132").Value = -5607.572
$ws.Range("H136").Value = 114618.89
$ws.Range("I136").Value = 3716.4285
$ws.Range("K136").Value = 11149.2855
$ws.Range("M136").Value = -8599.2855
